$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C ("Förändrad") for rows 2-13 from 45175 to 45183
for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 3).Value = 45183
}
